# edit.ps1 -- apply the "dissertation_plan" commit:
#   * append a new bullet to the "Minor problems" slide about the GCF/GCA mess
#   * insert a new "More problems (?)" slide (with an empty body) right after it
#
# Corresponds to the commit message:
#   "big step host stuff
#    did a lot of notebooking today, some more tables, about 80-90% there now"

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. "Minor problems" slide (position 10): add a trailing paragraph.
# ---------------------------------------------------------------------------
$minorProblems = $p.Slides.Item(10)
$body = $minorProblems.Shapes.Item(2)
$bodyRange = $body.TextFrame.TextRange

$newPara = $bodyRange.InsertAfter([char]13 + "The whole GCF / GCA ")
$newPara.IndentLevel = 0
$bodyRange.InsertAfter("debaucle") | Out-Null
$bodyRange.InsertAfter(" ") | Out-Null

# ---------------------------------------------------------------------------
# 2. Insert a brand new "Title and Content" slide right after it (position 11).
# ---------------------------------------------------------------------------
$titleAndContentLayout = $p.SlideMaster.CustomLayouts.Item(2)
$newSlide = $p.Slides.AddSlide(11, $titleAndContentLayout)
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "More problems (?)"
